# Delete single entry: "331 - CF" (08-12-2025 -> 08-15-2025), which lives in
# row 116 (Unique ID 114). Column A holds a static sequential index and must
# stay untouched, so instead of a normal row delete we shift only columns
# B:J up by one row for every row from the deleted entry through the end of
# the table, then drop the now-duplicated trailing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 116
$lastRow = 208

for ($r = $firstRow; $r -le ($lastRow - 1); $r++) {
    $srcRow = $r + 1
    $srcRange = "B$srcRow" + ":J$srcRow"
    $dstRange = "B$r" + ":J$r"
    $val = $ws.Range($srcRange).Value()
    $ws.Range($dstRange).Value = $val
}

$ws.Rows($lastRow).Delete()
